$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value while forcing Text storage so numeric-looking
# strings (e.g. "148.10", "28.00") keep their exact literal formatting
# instead of being auto-coerced into numbers, then restore the default
# "Normal" style so no stray NumberFormat style is left on the cell.
function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue 'D2' '61.573.30'
Set-TextValue 'E2' '  +1.28%  '
Set-TextValue 'D3' '3.450.31'
Set-TextValue 'E3' '  +2.40%  '
Set-TextValue 'E4' '  +0.02%  '
Set-TextValue 'D5' '579.39'
Set-TextValue 'E5' '  +1.41%  '
Set-TextValue 'D6' '148.10'
Set-TextValue 'E6' '  +8.85%  '
Set-TextValue 'D7' '3.452.27'
Set-TextValue 'E7' '  +2.54%  '
Set-TextValue 'E8' '  +0.09%  '
Set-TextValue 'E9' '  +0.71%  '
Set-TextValue 'D10' '7.77'
Set-TextValue 'E10' '  +3.48%  '
Set-TextValue 'E11' '  +0.82%  '
Set-TextValue 'D12' '0.392'
Set-TextValue 'E12' '  +1.62%  '
Set-TextValue 'D13' '4.041.96'
Set-TextValue 'E13' '  +2.65%  '
Set-TextValue 'D14' '28.00'
Set-TextValue 'E14' '  +6.78%  '
Set-TextValue 'E15' '  -0.70%  '
Set-TextValue 'E16' '  +1.48%  '
Set-TextValue 'D17' '3.446.50'
Set-TextValue 'E17' '  +2.37%  '
Set-TextValue 'D18' '61.705.05'
Set-TextValue 'E18' '  +1.25%  '
Set-TextValue 'E19' '  +8.70%  '
Set-TextValue 'D20' '14.32'
Set-TextValue 'E20' '  +2.13%  '
Set-TextValue 'D21' '9.43'
Set-TextValue 'E21' '  +1.36%  '
Set-TextValue 'E22' '  +2.72%  '
Set-TextValue 'D23' '0.569'
Set-TextValue 'E23' '  +2.54%  '
Set-TextValue 'D24' '3.597.93'
Set-TextValue 'E24' '  +2.78%  '
Set-TextValue 'B25' 'Litecoin'
Set-TextValue 'C25' 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue 'D25' '72.67'
Set-TextValue 'E25' '  +2.26%  '
Set-TextValue 'D26' '5.78'
Set-TextValue 'E26' '  +1.07%  '
Set-TextValue 'B27' 'Dai'
Set-TextValue 'C27' 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue 'D27' '0.999'
Set-TextValue 'E27' '  -0.17%  '
Set-TextValue 'E28' '  -1.61%  '
Set-TextValue 'E29' '  +7.60%  '
Set-TextValue 'D30' '7.82'
Set-TextValue 'E30' '  +4.41%  '
Set-TextValue 'D31' '0.999'
Set-TextValue 'E31' '  -0.09%  '
Set-TextValue 'E32' '  -14.26%  '
Set-TextValue 'E33' '  +1.47%  '
Set-TextValue 'E34' '  +1.46%  '
Set-TextValue 'E35' '  +0.01%  '
Set-TextValue 'D36' '23.98'
Set-TextValue 'E36' '  +1.22%  '
Set-TextValue 'E37' '  +4.19%  '
Set-TextValue 'D38' '5.22'
Set-TextValue 'E38' '  +0.27%  '
Set-TextValue 'E39' '  +2.56%  '
Set-TextValue 'D40' '166.23'
Set-TextValue 'E40' '  +0.83%  '
Set-TextValue 'E41' '  +4.91%  '
Set-TextValue 'D42' '26.09'
Set-TextValue 'E42' '  +10.18%  '
Set-TextValue 'E43' '  +3.61%  '
Set-TextValue 'E44' '  +0.12%  '
Set-TextValue 'D45' '4.50'
Set-TextValue 'E45' '  +2.32%  '
Set-TextValue 'D46' '42.34'
Set-TextValue 'E46' '  +2.07%  '
Set-TextValue 'E47' '  +1.48%  '
Set-TextValue 'B48' 'Maker'
Set-TextValue 'C48' 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue 'D48' '2.595.87'
Set-TextValue 'E48' '  +10.26%  '
Set-TextValue 'B49' 'ONDO'
Set-TextValue 'C49' 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
Set-TextValue 'D49' '1.16'
Set-TextValue 'E49' '  -3.13%  '
Set-TextValue 'D50' '6.97'
Set-TextValue 'E50' '  +2.56%  '
Set-TextValue 'D51' '23.30'
Set-TextValue 'E51' '  -0.01%  '
